$wb = $excel.ActiveWorkbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$openList = $wb.Worksheets.Add($null, $last)
$openList.Name = "OpenList"
$last2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$closedList = $wb.Worksheets.Add($null, $last2)
$closedList.Name = "ClosedList"
$openList.Activate()
Write-Host $wb.Worksheets.Count
